$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.834.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.47%  "
$ws.Range("D3").Value = "'1.694.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.77%  "
$ws.Range("D4").Value = "'0.9985"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'221.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.81%  "
$ws.Range("D6").Value = "'0.5141"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.93%  "
$ws.Range("D7").Value = "'0.9970"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").Value = "'0.2585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.99%  "
$ws.Range("D9").Value = "'22.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("D10").Value = "'0.06217"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.22%  "
$ws.Range("D11").Value = "'0.07337"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").Value = "'1.653.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.86%  "
$ws.Range("D13").Value = "'4.491"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.93%  "
$ws.Range("D14").Value = "'0.5823"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.60%  "
$ws.Range("D15").Value = "'1.879.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.98%  "
$ws.Range("D16").Value = "'0.000008149"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -13.45%  "
$ws.Range("D17").Value = "'65.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -12.32%  "
$ws.Range("D18").Value = "'26.740.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.06%  "
$ws.Range("D19").Value = "'5.014"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.62%  "
$ws.Range("D20").Value = "'1.007"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'10.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.63%  "
$ws.Range("D22").Value = "'183.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.89%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'6.274"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.33%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'143.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.41%  "
$ws.Range("D26").Value = "'7.601"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("D27").Value = "'0.1148"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.51%  "
$ws.Range("D28").Value = "'15.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.40%  "
$ws.Range("D29").Value = "'1.335"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.70%  "
$ws.Range("D30").Value = "'0.05885"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.83%  "
$ws.Range("D31").Value = "'1.346"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.77%  "
$ws.Range("D32").Value = "'3.451"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.14%  "
$ws.Range("D33").Value = "'3.417"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.87%  "
$ws.Range("D34").Value = "'1.658"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").Value = "'0.9927"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.33%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6032"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.22%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("D38").Value = "'2.638"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.92%  "
$ws.Range("D39").Value = "'0.01594"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.05%  "
$ws.Range("D40").Value = "'1.085.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("D41").Value = "'0.8589"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "'1.007"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "'5.765"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.15%  "
$ws.Range("D44").Value = "'97.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").Value = "'1.791.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.23%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'56.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.32%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.010"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").Value = "'0.00000000104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.09%  "
$ws.Range("D49").Value = "'0.4382"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").Value = "'0.05225"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.461"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.01%  "
